$d = $word.ActiveDocument

# 1) ISP section: "It's nicely explained on this website " -> "It's nicely explained in here "
$d.Content.Find.Execute(
    "nicely explained on this website ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "nicely explained in here ",
    2)

# 2) ISP section: "It's contra productive to rewrite it in this file :)"
#    -> "It would be contra productive to rewrite it into this file :)"
$d.Content.Find.Execute(
    "It’s contra productive to rewrite it in this file :)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "It would be contra productive to rewrite it into this file :)",
    2)

# 3) DIP solution: remove "which represents the logging"
$d.Content.Find.Execute(
    "an interface which represents the logging, this interface",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "an interface, this interface",
    2)

# 4) DIP solution: "parent term" -> "main term"
$d.Content.Find.Execute(
    "is a parent term while",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "is a main term while",
    2)

# 5) DIP solution: normalize the non-breaking space before "concept" to a
#    regular space (the sentence was retyped) and drop the trailing
#    "(Source: Google)" parenthetical
$d.Content.Find.Execute(
    "is also a" + [char]0x00A0 + "concept where the flow of the application is inverted (Source: Google)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "is also a concept where the flow of the application is inverted.",
    2)
